$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number (e.g. "1.004") must be forced
# back to Text so Excel does not auto-convert them to a numeric value - the
# source data stores prices/volumes as text, matching the original file layout.
$textNumberCells = @("D4", "D5", "D6", "D8", "D9", "D10", "D11", "D13", "D17", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D42", "D43", "D44", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textNumberCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "25.849.98"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.637.27"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "215.44"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "0.5089"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "0.2588"
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").Value = "0.06437"
$ws.Range("E9").Value = "  +1.93%  "
$ws.Range("D10").Value = "20.41"
$ws.Range("E10").Value = "  +5.25%  "
$ws.Range("D11").Value = "0.07804"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.659.76"
$ws.Range("E12").Value = "  +1.97%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.262"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").Value = "1.863.32"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").Value = "0.0₅7676"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").Value = "63.37"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "25.876.07"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("D19").Value = "1.004"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "193.36"
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "4.386"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").Value = "9.973"
$ws.Range("E22").Value = "  +1.48%  "
$ws.Range("D23").Value = "6.161"
$ws.Range("E23").Value = "  +3.02%  "
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "1.755"
$ws.Range("E25").Value = "  -6.95%  "
$ws.Range("D26").Value = "139.08"
$ws.Range("E26").Value = "  -2.07%  "
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").Value = "6.847"
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("D29").Value = "15.56"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("E30").Value = "  +0.29%  "
$ws.Range("D31").Value = "0.04978"
$ws.Range("E31").Value = "  +2.20%  "
$ws.Range("D32").Value = "3.309"
$ws.Range("E32").Value = "  +2.73%  "
$ws.Range("D33").Value = "3.258"
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("D34").Value = "1.571"
$ws.Range("D35").Value = "2.387"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").Value = "0.9039"
$ws.Range("E36").Value = "  +1.47%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").Value = "2.577"
$ws.Range("E37").Value = "  +1.68%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.5584"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("D39").Value = "1.134.19"
$ws.Range("E39").Value = "  +2.21%  "
$ws.Range("D40").Value = "0.01573"
$ws.Range("E40").Value = "  +1.84%  "
$ws.Range("E41").Value = "  -0.33%  "
$ws.Range("D42").Value = "5.467"
$ws.Range("E42").Value = "  -1.51%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "99.21"
$ws.Range("E43").Value = "  +2.14%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "0.8000"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("D45").Value = "0.0₈113"
$ws.Range("E45").Value = "  +1.47%  "
$ws.Range("D46").Value = "55.73"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("D47").Value = "0.4268"
$ws.Range("E47").Value = "  -3.66%  "
$ws.Range("D48").Value = "7.832"
$ws.Range("E48").Value = "  +4.28%  "
$ws.Range("D49").Value = "0.05030"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").Value = "1.000"
$ws.Range("E50").Value = "  +0.62%  "
$ws.Range("E51").Value = "  +0.52%  "
